$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 9 (Louisville City FC vs Detroit City FC)
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Louisville City FC  - Detroit City FC: 22:00"
$ws.Range("B9").Value = "Louisville City FC"
$ws.Range("C9").Value = 64
$ws.Range("D9").Value = 75
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 1.53
$ws.Range("G9").Value = ""

# Insert two new rows before current row 12 (after the shift caused by the row above)
$ws.Range("12:13").Insert()

$ws.Range("A12").Value = "Kuala Lumpur City FC ✓ - Kelantan The Real Warriors: 2:1"
$ws.Range("B12").Value = "Kuala Lumpur City FC"
$ws.Range("C12").Value = 62
$ws.Range("D12").Value = 75
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 1.53
$ws.Range("G12").Value = "✓"

$ws.Range("A13").Value = "Thimphu City FC  - BFF Academy FC: -:-'"
$ws.Range("B13").Value = "Thimphu City FC"
$ws.Range("C13").Value = 62
$ws.Range("D13").Value = 75
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 1.53
$ws.Range("G13").Value = ""

# Append three new rows at the end (rows 32-34)
$ws.Range("A32").Value = "Leicester City X - Blackburn Rovers: 0:2"
$ws.Range("B32").Value = "Leicester City"
$ws.Range("C32").Value = 28
$ws.Range("D32").Value = 75
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = 1.53
$ws.Range("G32").Value = "X"

$ws.Range("A33").Value = "Cork City FC - Derry City : -:-'"
$ws.Range("B33").Value = "Derry City"
$ws.Range("C33").Value = 26
$ws.Range("D33").Value = 75
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = 1.53
$ws.Range("G33").Value = ""

$ws.Range("A34").Value = "PDRM FC - Kuching City ✓: 0:5"
$ws.Range("B34").Value = "Kuching City"
$ws.Range("C34").Value = 25
$ws.Range("D34").Value = 75
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = 1.53
$ws.Range("G34").Value = "✓"
